$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$chartObj = $ws.ChartObjects().Add(542520, 542520, 4459680, 4459680)
$chart = $chartObj.Chart
$chart.ChartType = 74
$chart.SetSourceData($ws.Range("A1:B10"))
$s2 = $chart.SeriesCollection().NewSeries()
$s2.AxisGroup = 2

$ax_p = $chart.Axes(2,1)
$ax_p.HasTitle = $true
$ax_p.AxisTitle.Text = "PRIMARY_AXIS_TITLE_UNIQUE"

$ax_s = $chart.Axes(2,2)
$ax_s.HasTitle = $true
$ax_s.AxisTitle.Text = "SECONDARY_AXIS_TITLE_UNIQUE"

Write-Host "Readback primary: $($ax_p.AxisTitle.Text)"
Write-Host "Readback secondary: $($ax_s.AxisTitle.Text)"
